$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$src = $ws.Range("A9:D29")
$dst = $ws.Range("A10:D30")
$src.Copy()
$dst.PasteSpecial(-4122)
$dst.Value2 = $src.Value2
$excel.CutCopyMode = 0
Write-Output "done"
